# Applies the "Add implementation for storage management in DG" edit to
# docs/diagrams/ExportSequenceDiagram.pptx (slide 1):
#   1. Extend the ":LogicManager" lifeline connector downward.
#   2. Extend the ":Storage" lifeline connector downward.
#   3. Rename the "exportDeskBoard(" call label to "saveDeskBoard(" while
#      keeping the surrounding "deskBoard, filePath.value)" text intact -
#      this is modeled as a new leading run "saveDeskBoard" plus the
#      existing run's text shrinking down to just "(".

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# EMU -> Points helper. PowerPoint COM shape geometry (Top/Left/Width/
# Height) is expressed in points (1 pt = 12700 EMU) and stored/returned
# as a 32-bit float, which can round a straight "/12700.0" down by a
# single EMU on save. Nudging by half an EMU before dividing keeps the
# round-trip exact.
function EmuToPoints($emu) {
    return ($emu + 0.5) / 12700.0
}

# --- 1) ":LogicManager" lifeline connector ("Straight Connector 4") ---
$lifelineLogicManager = $s.Shapes.Item(4)
$lifelineLogicManager.Height = EmuToPoints 3975729

# --- 2) ":Storage" lifeline connector ("Straight Connector 93") ---
$lifelineStorage = $s.Shapes.Item(49)
$lifelineStorage.Height = EmuToPoints 3091970

# --- 3) "exportDeskBoard(...)" call label -> "saveDeskBoard(...)" ---
$callLabel = $s.Shapes.Item(46)
$tr = $callLabel.TextFrame.TextRange

# First shrink the original "exportDeskBoard(" run down to "(" - a
# same-length-reduction in place (no InsertBefore) so the textbox's
# autosized height isn't perturbed.
$openParenRun = $tr.Characters(1, 16)
$openParenRun.Text = "saveDeskBoard("

# Now split that merged run back into two runs - "saveDeskBoard" and
# "(" - by re-asserting the first 13 characters; this is a same-length
# no-op textually but forces the run boundary PowerPoint needs so the
# two pieces keep/gain independent formatting.
$newNameRun = $tr.Characters(1, 13)
$newNameRun.Text = "saveDeskBoard"
